$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Razon social" / "Nombre Fantasia" text fixes: commas used to separate
#    multiple co-holders on a contract were replaced with periods (and a
#    couple of "S.H." abbreviations were tightened to "SH").
# ---------------------------------------------------------------------------
$nameFixes = @(
    @("E28",  "TRABICHET MARIA. VERGARA ADEL Y OTRA"),
    @("F28",  "TRABICHET MARIA. VERGARA ADEL Y OTRA"),
    @("E125", "TRABICHET MARIA. VERGARA ADEL Y OTRA"),
    @("F125", "TRABICHET MARIA. VERGARA ADEL Y OTRA"),

    @("E36",  "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"),
    @("F36",  "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"),
    @("E93",  "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"),
    @("F93",  "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"),

    @("E42",  "RAMIREZ CLAUDIA. RAMIREZ CESAR Y RAMIREZ VERONICA SH"),

    @("E43",  "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"),
    @("E209", "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"),

    @("E50",  "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"),
    @("E192", "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"),

    @("E94",  "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"),

    @("E140", "RICCOTTI. MARIANA EDITH"),

    @("E183", "FERNANDEZ. MARIO HUGO")
)

foreach ($fix in $nameFixes) {
    $ws.Range($fix[0]).Value = $fix[1]
}

# ---------------------------------------------------------------------------
# 2) "Importe" column (H2:H249) was scraped with Spanish/Argentine number
#    formatting (periods as thousands separators, commas as decimal
#    separators) but stored as plain text. Re-write every value as a plain
#    dot-decimal numeric string (e.g. "2.040,00" -> "2040.00"), keeping the
#    values as text.
# ---------------------------------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 249 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = [string]$cell.Value2
    if ([string]::IsNullOrEmpty($old)) { continue }
    $new = $old.Replace(".", "").Replace(",", ".")
    $cell.NumberFormat = "@"
    $cell.Value = $new
}
